# The deck's applied design ("Integral") and the stock "Office Theme"
# that ships alongside it (previously only reachable via the Notes
# Master) swap places: the theme actually driving the slides / slide
# master switches from the Integral palette to the plain Office Theme
# palette.
#
# Mirror what a user does from Slide Master view -> Colors -> Customize
# Colors...: walk the twelve theme colour slots (Text/Background Dark
# 1, Light 1, Dark 2, Light 2, Accent 1-6, Hyperlink, Followed
# Hyperlink) on the presentation's theme and set each one to the Office
# Theme's RGB value.
#
# (PowerPoint's ThemeColor.RGB takes a COM BGR-packed long -
# R + G*256 + B*65536 - hence the decimal values below instead of the
# more familiar RRGGBB hex form.)

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
